# Refresh market-price-derived profit figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR leve-profit sheets, as
# produced by the scheduled market-data runner. CUL has no changes this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3149.139
$ws.Range("I64").Value = 3067.8333
$ws.Range("J64").Value = 3165.4
$ws.Range("K64").Value = 3067.8333
$ws.Range("L64").Value = 3165.4
$ws.Range("M64").Value = -2819.8333
$ws.Range("N64").Value = -3661.4

$ws.Range("H67").Value = 3149.139
$ws.Range("I67").Value = 3067.8333
$ws.Range("J67").Value = 3165.4
$ws.Range("K67").Value = 3067.8333
$ws.Range("L67").Value = 3165.4
$ws.Range("M67").Value = -2209.8333
$ws.Range("N67").Value = -4881.4

$ws.Range("H74").Value = 4227.727
$ws.Range("I74").Value = 4523.846
$ws.Range("K74").Value = 4523.846
$ws.Range("M74").Value = -3587.846

$ws.Range("H77").Value = 4227.727
$ws.Range("I77").Value = 4523.846
$ws.Range("K77").Value = 22619.23
$ws.Range("M77").Value = -17939.23

$ws.Range("H86").Value = 7923.1665
$ws.Range("I86").Value = 14312.875
$ws.Range("J86").Value = 2811.4
$ws.Range("K86").Value = 14312.875
$ws.Range("L86").Value = 2811.4
$ws.Range("M86").Value = -13189.875
$ws.Range("N86").Value = -5057.4

$ws.Range("H89").Value = 7923.1665
$ws.Range("I89").Value = 14312.875
$ws.Range("J89").Value = 2811.4
$ws.Range("K89").Value = 71564.375
$ws.Range("L89").Value = 14057
$ws.Range("M89").Value = -65948.375
$ws.Range("N89").Value = -25289

$ws.Range("H137").Value = 1629.1724
$ws.Range("I137").Value = 1221.75
$ws.Range("K137").Value = 3665.25
$ws.Range("M137").Value = -1115.25


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 440.33334
$ws.Range("I2").Value = 440.33334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 440.33334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -327.33334
$ws.Range("N2").ClearContents()

$ws.Range("H61").Value = 2094.3076
$ws.Range("I61").Value = 1825.1111
$ws.Range("J61").Value = 2700
$ws.Range("K61").Value = 1825.1111
$ws.Range("L61").Value = 2700
$ws.Range("M61").Value = -1613.1111
$ws.Range("N61").Value = -3124

$ws.Range("H63").Value = 4235.25
$ws.Range("I63").Value = 2474
$ws.Range("J63").Value = 5292
$ws.Range("K63").Value = 2474
$ws.Range("L63").Value = 5292
$ws.Range("M63").Value = -1788
$ws.Range("N63").Value = -6664

$ws.Range("H66").Value = 4235.25
$ws.Range("I66").Value = 2474
$ws.Range("J66").Value = 5292
$ws.Range("K66").Value = 12370
$ws.Range("L66").Value = 26460
$ws.Range("M66").Value = -8938
$ws.Range("N66").Value = -33324

$ws.Range("H88").Value = 58825732
$ws.Range("I88").Value = 1663.125
$ws.Range("J88").Value = 111113790
$ws.Range("K88").Value = 1663.125
$ws.Range("L88").Value = 111113790
$ws.Range("M88").Value = -1257.125
$ws.Range("N88").Value = -111114602

$ws.Range("H91").Value = 58825732
$ws.Range("I91").Value = 1663.125
$ws.Range("J91").Value = 111113790
$ws.Range("K91").Value = 1663.125
$ws.Range("L91").Value = 111113790
$ws.Range("M91").Value = -259.125
$ws.Range("N91").Value = -111116598

$ws.Range("H116").Value = 440.33334
$ws.Range("I116").Value = 440.33334
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 440.33334
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1853.66666
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 2781.4119
$ws.Range("I122").Value = 2498.8572
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 7496.571599999999
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -5046.571599999999
$ws.Range("N122").Value = -17200

$ws.Range("H132").Value = 1830.0857
$ws.Range("I132").Value = 1513.1154
$ws.Range("J132").Value = 2745.7778
$ws.Range("K132").Value = 4539.3462
$ws.Range("L132").Value = 8237.3334
$ws.Range("M132").Value = -2009.3462
$ws.Range("N132").Value = -13297.3334

$ws.Range("H136").Value = 2094.3076
$ws.Range("I136").Value = 1825.1111
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 5475.3333
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -2925.3333
$ws.Range("N136").Value = -13200


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 440.33334
$ws.Range("I3").Value = 440.33334
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 440.33334
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -326.33334
$ws.Range("N3").ClearContents()

$ws.Range("H86").Value = 40002760
$ws.Range("I86").Value = 50002076
$ws.Range("K86").Value = 50002076
$ws.Range("M86").Value = -50000953

$ws.Range("H89").Value = 40002760
$ws.Range("I89").Value = 50002076
$ws.Range("K89").Value = 250010380
$ws.Range("M89").Value = -250004764


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2262.4211
$ws.Range("I132").Value = 1558.2222
$ws.Range("J132").Value = 2896.2
$ws.Range("K132").Value = 4674.6666
$ws.Range("L132").Value = 8688.599999999999
$ws.Range("M132").Value = -2144.6666
$ws.Range("N132").Value = -13748.6


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1523.85
$ws.Range("I113").Value = 1116.5454
$ws.Range("J113").Value = 2021.6666
$ws.Range("K113").Value = 1116.5454
$ws.Range("L113").Value = 2021.6666
$ws.Range("M113").Value = 1053.4546
$ws.Range("N113").Value = -6361.6666

$ws.Range("H122").Value = 3461.361
$ws.Range("I122").Value = 2949.75
$ws.Range("J122").Value = 4484.5835
$ws.Range("K122").Value = 8849.25
$ws.Range("L122").Value = 13453.7505
$ws.Range("M122").Value = -6399.25
$ws.Range("N122").Value = -18353.7505

$ws.Range("H132").Value = 2697.1052
$ws.Range("I132").Value = 2151.5833
$ws.Range("J132").Value = 3632.2856
$ws.Range("K132").Value = 6454.749899999999
$ws.Range("L132").Value = 10896.8568
$ws.Range("M132").Value = -3924.749899999999
$ws.Range("N132").Value = -15956.8568


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2554699.5
$ws.Range("I136").Value = 4813347.5
$ws.Range("K136").Value = 14440042.5
$ws.Range("M136").Value = -14437492.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 18984.615
$ws.Range("J64").Value = 18984.615
$ws.Range("L64").Value = 18984.615
$ws.Range("N64").Value = -19480.615

$ws.Range("H67").Value = 18984.615
$ws.Range("J67").Value = 18984.615
$ws.Range("L67").Value = 18984.615
$ws.Range("N67").Value = -20700.615

